$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H17").Value = 1118.5814
$ws_ALC.Range("J17").Value = 1118.5814
$ws_ALC.Range("L17").Value = 3355.7442
$ws_ALC.Range("N17").Value = -3691.7442

$ws_ALC.Range("H100").Value = 4217.5386
$ws_ALC.Range("I100").Value = 2475
$ws_ALC.Range("J100").Value = 4992
$ws_ALC.Range("K100").Value = 2475
$ws_ALC.Range("L100").Value = 4992
$ws_ALC.Range("M100").Value = -1934
$ws_ALC.Range("N100").Value = -6074

$ws_ALC.Range("H103").Value = 484.6111
$ws_ALC.Range("J103").Value = 478
$ws_ALC.Range("L103").Value = 1434
$ws_ALC.Range("N103").Value = -2606

$ws_ALC.Range("H111").Value = 2430.3333
$ws_ALC.Range("I111").Value = 1598.25
$ws_ALC.Range("J111").Value = 3096
$ws_ALC.Range("K111").Value = 4794.75
$ws_ALC.Range("L111").Value = 9288
$ws_ALC.Range("M111").Value = -1727.75
$ws_ALC.Range("N111").Value = -15422

$ws_ALC.Range("H112").Value = 1857.88
$ws_ALC.Range("J112").Value = 1948.174
$ws_ALC.Range("L112").Value = 5844.522
$ws_ALC.Range("N112").Value = -8060.522

$ws_ALC.Range("H116").Value = 4267.4614
$ws_ALC.Range("I116").Value = 4276.4443
$ws_ALC.Range("J116").Value = 4247.25
$ws_ALC.Range("K116").Value = 4276.4443
$ws_ALC.Range("L116").Value = 4247.25
$ws_ALC.Range("M116").Value = -834.4443000000001
$ws_ALC.Range("N116").Value = -11131.25

$ws_ALC.Range("H132").Value = 972.16
$ws_ALC.Range("J132").Value = 815.75
$ws_ALC.Range("L132").Value = 2447.25
$ws_ALC.Range("N132").Value = -7507.25

$ws_ALC.Range("H135").Value = 1646.25
$ws_ALC.Range("I135").Value = 1651.8518
$ws_ALC.Range("K135").Value = 14866.6662
$ws_ALC.Range("M135").Value = -12331.6662

$ws_ALC.Range("H137").Value = 1926.5
$ws_ALC.Range("J137").Value = 1323.5264
$ws_ALC.Range("L137").Value = 3970.5792
$ws_ALC.Range("N137").Value = -9070.5792

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H61").Value = 1909.3572
$ws_ARM.Range("I61").Value = 1848.5454
$ws_ARM.Range("J61").Value = 2132.3333
$ws_ARM.Range("K61").Value = 1848.5454
$ws_ARM.Range("L61").Value = 2132.3333
$ws_ARM.Range("M61").Value = -1636.5454
$ws_ARM.Range("N61").Value = -2556.3333

$ws_ARM.Range("H74").Value = 2989
$ws_ARM.Range("I74").Value = 2657.4285
$ws_ARM.Range("J74").Value = 3569.25
$ws_ARM.Range("K74").Value = 2657.4285
$ws_ARM.Range("L74").Value = 3569.25
$ws_ARM.Range("M74").Value = -1783.4285
$ws_ARM.Range("N74").Value = -5317.25

$ws_ARM.Range("H77").Value = 2989
$ws_ARM.Range("I77").Value = 2657.4285
$ws_ARM.Range("J77").Value = 3569.25
$ws_ARM.Range("K77").Value = 13287.1425
$ws_ARM.Range("L77").Value = 17846.25
$ws_ARM.Range("M77").Value = -8919.1425
$ws_ARM.Range("N77").Value = -26582.25

$ws_ARM.Range("H122").Value = 4640.548
$ws_ARM.Range("J122").Value = 4874.375
$ws_ARM.Range("L122").Value = 14623.125
$ws_ARM.Range("N122").Value = -19523.125

$ws_ARM.Range("H132").Value = 7992.7446
$ws_ARM.Range("I132").Value = 4186.35
$ws_ARM.Range("K132").Value = 12559.05
$ws_ARM.Range("M132").Value = -10029.05

$ws_ARM.Range("H136").Value = 1909.3572
$ws_ARM.Range("I136").Value = 1848.5454
$ws_ARM.Range("J136").Value = 2132.3333
$ws_ARM.Range("K136").Value = 5545.6362
$ws_ARM.Range("L136").Value = 6396.999899999999
$ws_ARM.Range("M136").Value = -2995.6362
$ws_ARM.Range("N136").Value = -11496.9999

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H134").Value = 3125.111
$ws_BSM.Range("I134").Value = 2650.6206
$ws_BSM.Range("K134").Value = 7951.861800000001
$ws_BSM.Range("M134").Value = -5416.861800000001

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 1633.125
$ws_CRP.Range("I31").Value = 1466.3334
$ws_CRP.Range("K31").Value = 1466.3334
$ws_CRP.Range("M31").Value = -1171.3334

$ws_CRP.Range("H34").Value = 1633.125
$ws_CRP.Range("I34").Value = 1466.3334
$ws_CRP.Range("K34").Value = 1466.3334
$ws_CRP.Range("M34").Value = -1264.3334

$ws_CRP.Range("H58").Value = 5556.625
$ws_CRP.Range("I58").Value = 7310.8
$ws_CRP.Range("K58").Value = 7310.8
$ws_CRP.Range("M58").Value = -7107.8

$ws_CRP.Range("H134").Value = 1578.6923
$ws_CRP.Range("I134").Value = 1415.4286
$ws_CRP.Range("J134").Value = 1769.1666
$ws_CRP.Range("K134").Value = 4246.2858
$ws_CRP.Range("L134").Value = 5307.4998
$ws_CRP.Range("M134").Value = -1711.2858
$ws_CRP.Range("N134").Value = -10377.4998

$ws_CRP.Range("H136").Value = 5556.625
$ws_CRP.Range("I136").Value = 7310.8
$ws_CRP.Range("K136").Value = 21932.4
$ws_CRP.Range("M136").Value = -19382.4

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H82").Value = 13179.75
$ws_CUL.Range("I82").Value = 10906.333
$ws_CUL.Range("J82").Value = 20000
$ws_CUL.Range("K82").Value = 32718.999
$ws_CUL.Range("L82").Value = 60000
$ws_CUL.Range("M82").Value = -32312.999
$ws_CUL.Range("N82").Value = -60812

$ws_CUL.Range("H85").Value = 13179.75
$ws_CUL.Range("I85").Value = 10906.333
$ws_CUL.Range("J85").Value = 20000
$ws_CUL.Range("K85").Value = 32718.999
$ws_CUL.Range("L85").Value = 60000
$ws_CUL.Range("M85").Value = -31314.999
$ws_CUL.Range("N85").Value = -62808

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H122").Value = 2459.625
$ws_GSM.Range("I122").Value = 2006.0416
$ws_GSM.Range("J122").Value = 3820.375
$ws_GSM.Range("K122").Value = 6018.1248
$ws_GSM.Range("L122").Value = 11461.125
$ws_GSM.Range("M122").Value = -3568.1248
$ws_GSM.Range("N122").Value = -16361.125

$ws_GSM.Range("H132").Value = 2401.3333
$ws_GSM.Range("I132").Value = 1990.7693
$ws_GSM.Range("K132").Value = 5972.3079
$ws_GSM.Range("M132").Value = -3442.3079

$ws_GSM.Range("H134").Value = 80000
$ws_GSM.Range("J134").Value = 80000
$ws_GSM.Range("L134").Value = 240000
$ws_GSM.Range("N134").Value = -245070

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H16").Value = 1250.64
$ws_LTW.Range("I16").Value = 1044.9333
$ws_LTW.Range("K16").Value = 1044.9333
$ws_LTW.Range("M16").Value = -874.9332999999999

$ws_LTW.Range("H22").Value = 2664.4092
$ws_LTW.Range("I22").Value = 1030
$ws_LTW.Range("K22").Value = 1030
$ws_LTW.Range("M22").Value = -735

$ws_LTW.Range("H27").Value = 2664.4092
$ws_LTW.Range("I27").Value = 1030
$ws_LTW.Range("K27").Value = 1030
$ws_LTW.Range("M27").Value = -923

$ws_LTW.Range("H46").Value = 1746.8
$ws_LTW.Range("I46").Value = 1442
$ws_LTW.Range("J46").Value = 1950
$ws_LTW.Range("K46").Value = 1442
$ws_LTW.Range("L46").Value = 1950
$ws_LTW.Range("M46").Value = -1254
$ws_LTW.Range("N46").Value = -2326

$ws_LTW.Range("H136").Value = 2443.6191
$ws_LTW.Range("I136").Value = 2196.889
$ws_LTW.Range("K136").Value = 6590.667
$ws_LTW.Range("M136").Value = -4040.667

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H74").Value = 30446
$ws_WVR.Range("J74").Value = 30446
$ws_WVR.Range("L74").Value = 30446
$ws_WVR.Range("N74").Value = -32318

$ws_WVR.Range("H77").Value = 30446
$ws_WVR.Range("J77").Value = 30446
$ws_WVR.Range("L77").Value = 91338
$ws_WVR.Range("N77").Value = -100698

$ws_WVR.Range("H100").Value = 2993
$ws_WVR.Range("I100").Value = 989.5
$ws_WVR.Range("K100").Value = 1979
$ws_WVR.Range("M100").Value = -1438

$ws_WVR.Range("H132").Value = 1749.4
$ws_WVR.Range("I132").Value = 1699.4
$ws_WVR.Range("K132").Value = 5098.200000000001
$ws_WVR.Range("M132").Value = -2568.200000000001

$ws_WVR.Range("H136").Value = 5216.069
$ws_WVR.Range("I136").Value = 3452.577
$ws_WVR.Range("K136").Value = 10357.731
$ws_WVR.Range("M136").Value = -7807.731
